# Replace xlsx test files
# Updates the "code" column values for the charting question rows on the
# "Core" and "Test Chart" sheets of the import-test fixture, swapping the
# old generic placeholder codes for the new, more descriptive ones.

$wb = $excel.ActiveWorkbook
$core = $wb.Worksheets.Item("Core")
$chart = $wb.Worksheets.Item("Test Chart")

# -- Core sheet: rows 2-5, column A ("code") --
$core.Range("A2").Value = "ComplexChartInstanceName"
$core.Range("A3").Value = "ComplexChartDate"
$core.Range("A4").Value = "ComplexChartType"
$core.Range("A5").Value = "ComplexChartSubtype"

# -- Test Chart sheet: row 2, column A ("code") --
$chart.Range("A2").Value = "PatientChartingDate"

# These rows were re-entered as part of the fixture refresh, which is
# reflected in Excel re-assigning their vertical alignment (bottom, the
# default) explicitly instead of leaving it unset.
$core.Range("A2").VerticalAlignment = -4107
$core.Range("A3").VerticalAlignment = -4107
$core.Range("A4").VerticalAlignment = -4107
$core.Range("A5").VerticalAlignment = -4107
$chart.Range("A2").VerticalAlignment = -4107
